# Apply the data corrections recorded in the commit "Set up autopush for Saurav".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MiddleTRNew")

# A handful of outlier data points were retyped on the MiddleTRNew sheet.
$ws.Range("B9").Value = 431.139786666666
$ws.Range("A22").Value = 373.79067999999899
$ws.Range("C22").Value = 765.61044000000004
$ws.Range("E22").Value = 8047.6313600000003
$ws.Range("F22").Value = 1868.92659999999
$ws.Range("D23").Value = 3195.40466666666

# Reflect where the author's cursor ended up.
[void]$ws.Range("G18").Select()
